$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- User table (G/H columns): "string" -> "text" for existing rows ---
$ws.Range("H3").Value = "text"
$ws.Range("H4").Value = "text"
$ws.Range("H5").Value = "text"
$ws.Range("H6").Value = "text"

# --- Row 7 additions: Order.date and User.birth_date ---
$ws.Range("D7").Value = "date"
$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)

$ws.Range("E7").Value = "text"
$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122)

$ws.Range("G7").Value = "birth_date"
$ws.Range("G6").Copy()
$ws.Range("G7").PasteSpecial(-4122)

$ws.Range("H7").Value = "text"
$ws.Range("H6").Copy()
$ws.Range("H7").PasteSpecial(-4122)

# --- Row 8 additions: Product.date and User.start_date ---
$ws.Range("A8").Value = "date"
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B8").Value = "text"
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("G8").Value = "start_date"
$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)

$ws.Range("H8").Value = "text"
$ws.Range("H7").Copy()
$ws.Range("H8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- View state tweak: move the active selection like the authored edit did ---
$ws.Range("C15").Select()
